# Fruta / hortaliza, semanal
# Update the daily price records for "Macroferia Regional de Talca - Granada"
# and append a new weekly record (row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write one full data row (columns D, L, M, N, O, P, Q, R, S, T).
function Set-Row {
    param(
        [int]$r,
        [double]$d,
        [string]$l,
        [double]$m,
        [double]$n,
        [double]$o,
        [double]$p,
        [string]$q,
        [string]$rg,
        [double]$s,
        [double]$t
    )
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rg
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
}

Set-Row 2  44342 "Especial" 300 20000 20000 20000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1111 18
Set-Row 4  44299 "Primera"  100 15000 15000 15000 "`$/caja 15 kilos granel" "Provincia de Curicó" 1000 15
Set-Row 5  44354 "Primera"  100 18000 18000 18000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1000 18
Set-Row 6  44328 "Especial" 250 20000 20000 20000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1111 18
Set-Row 7  44348 "Especial" 200 20000 20000 20000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1111 18
Set-Row 8  44680 "Primera"  200 15000 15000 15000 "`$/caja 15 kilos granel" "Provincia de Limarí" 1000 15
Set-Row 9  44355 "Especial" 50  18000 18000 18000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1000 18
Set-Row 10 44326 "Especial" 300 20000 20000 20000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1111 18
Set-Row 11 44340 "Primera"  230 20000 20000 20000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1111 18
Set-Row 12 44319 "Especial" 120 20000 20000 20000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1111 18
Set-Row 13 44291 "Primera"  150 12000 12000 12000 "`$/caja 15 kilos granel" "Región Metropolitana" 800 15
Set-Row 14 44692 "Especial" 150 17000 17000 17000 "`$/caja 18 kilos granel" "Provincia de Limarí" 944 18
Set-Row 15 44294 "Primera"  50  12000 12000 12000 "`$/caja 15 kilos granel" "Región Metropolitana" 800 15
Set-Row 16 44316 "Especial" 300 20000 20000 20000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1111 18
Set-Row 17 44358 "Primera"  150 18000 18000 18000 "`$/caja 18 kilos granel" "Provincia de Limarí" 1000 18

# New row 18: same master product data as the rest of the sheet, with its own
# date / quality / volume / price / origin record.
$ws.Cells.Item(18, 1).Value = 5
$ws.Cells.Item(18, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(18, 3).Value = "Maule"
$ws.Cells.Item(18, 4).Value = 44358
$ws.Cells.Item(18, 4).NumberFormat = $ws.Range("D17").NumberFormat
$ws.Cells.Item(18, 5).Value = 7
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100104
$ws.Cells.Item(18, 8).Value = "Frutos de pepita"
$ws.Cells.Item(18, 9).Value = 100104001
$ws.Cells.Item(18, 10).Value = "Granada"
$ws.Cells.Item(18, 11).Value = "Wonderfull"
Set-Row 18 44358 "Primera" 100 17000 17000 17000 "`$/caja 18 kilos granel" "Provincia de Limarí" 944 18

$ws.Range("A1").Select()
